$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.446.02'
$ws.Range("E2").Value = '  +3.28%  '
$ws.Range("D3").Value = '2.082.89'
$ws.Range("E3").Value = '  +3.53%  '
$ws.Range("E4").Value = '  +0.21%  '
$ws.Range("D5").Value = '252.65'
$ws.Range("E5").Value = '  +2.42%  '
$ws.Range("D6").Value = '0.652'
$ws.Range("E6").Value = '  +1.12%  '
$ws.Range("D7").Value = '64.93'
$ws.Range("E7").Value = '  +2.32%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  +6.70%  '
$ws.Range("D10").Value = '59.70'
$ws.Range("E10").Value = '  +1.31%  '
$ws.Range("D11").Value = '0.0815'
$ws.Range("E11").Value = '  +8.61%  '
$ws.Range("E12").Value = '  -0.03%  '
$ws.Range("D13").Value = '0.932'
$ws.Range("E13").Value = '  -1.71%  '
$ws.Range("D14").Value = '23.51'
$ws.Range("E14").Value = '  +18.81%  '
$ws.Range("E15").Value = '  -1.29%  '
$ws.Range("D16").Value = '2.386.56'
$ws.Range("E16").Value = '  +3.62%  '
$ws.Range("D17").Value = '5.65'
$ws.Range("E17").Value = '  +3.50%  '
$ws.Range("D18").Value = '2.074.48'
$ws.Range("E18").Value = '  +2.95%  '
$ws.Range("D19").Value = '37.430.24'
$ws.Range("E19").Value = '  +3.48%  '
$ws.Range("D20").Value = '73.78'
$ws.Range("E20").Value = '  +2.13%  '
$ws.Range("D21").Value = '0.0₃0901'
$ws.Range("E21").Value = '  +4.69%  '
$ws.Range("D22").Value = '5.52'
$ws.Range("E22").Value = '  +4.71%  '
$ws.Range("D23").Value = '240.26'
$ws.Range("E23").Value = '  +2.43%  '
$ws.Range("E24").Value = '  -0.04%  '
$ws.Range("E25").Value = '  -1.91%  '
$ws.Range("D26").Value = '2.37'
$ws.Range("E26").Value = '  +3.14%  '
$ws.Range("D27").Value = '10.04'
$ws.Range("E27").Value = '  +3.01%  '
$ws.Range("D28").Value = '20.90'
$ws.Range("E28").Value = '  +6.04%  '
$ws.Range("D29").Value = '161.98'
$ws.Range("E29").Value = '  -2.61%  '
$ws.Range("D30").Value = '0.127'
$ws.Range("E30").Value = '  +27.95%  '
$ws.Range("E31").Value = '  +2.18%  '
$ws.Range("D32").Value = '5.18'
$ws.Range("E32").Value = '  +0.97%  '
$ws.Range("E33").Value = '  +3.50%  '
$ws.Range("D34").Value = '0.0629'
$ws.Range("E34").Value = '  +3.16%  '
$ws.Range("D35").Value = '4.68'
$ws.Range("E35").Value = '  +3.74%  '
$ws.Range("D36").Value = '2.55'
$ws.Range("E36").Value = '  +2.28%  '
$ws.Range("E37").Value = '  +11.31%  '
$ws.Range("E38").Value = '  +0.27%  '
$ws.Range("E39").Value = '  +2.58%  '
$ws.Range("D40").Value = '3.05'
$ws.Range("E40").Value = '  +29.65%  '
$ws.Range("B41").Value = 'Cronos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D41").Value = '0.103'
$ws.Range("E41").Value = '  +6.23%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = '1.29'
$ws.Range("E42").Value = '  +3.53%  '
$ws.Range("E43").Value = '  +4.41%  '
$ws.Range("D44").Value = '1.18'
$ws.Range("E44").Value = '  +4.69%  '
$ws.Range("E45").Value = '  +1.35%  '
$ws.Range("D46").Value = '17.25'
$ws.Range("E46").Value = '  +1.96%  '
$ws.Range("D47").Value = '96.18'
$ws.Range("E47").Value = '  +1.43%  '
$ws.Range("D48").Value = '8.00'
$ws.Range("E48").Value = '  +1.61%  '
$ws.Range("D49").Value = '1.403.30'
$ws.Range("E49").Value = '  +1.94%  '
$ws.Range("E50").Value = '  +0.98%  '
$ws.Range("D51").Value = '46.67'
$ws.Range("E51").Value = '  -1.31%  '
